$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "You are applying the K-nearest neighbors regressor and want the value of the new observation to be equal to the average value of all data points.  What should be the number of neighbors for the K-nearest neighbor regressor?",
        "ques_type": 2,
        "options": [
            "Number of current data points",
            "One",
            "Number of current data points minus one",
            "Number of current data points plus one"
        ],
        "score": "Number of current data points"
    },
    {
        "title": "You are applying gender classification and want to estimate the portion of correctly classified observation.  Which metric should you use?",
        "ques_type": 2,
        "options": [
            "Duality",
            "Accuracy",
            "Precision",
            "Recall"
        ],
        "score": "Accuracy"
    },
    {
        "title": "You are applying a decision-tree regressor on house prices and have got 12,345 as the value of the residual sum of squares (RSS).  How should you find an acceptable value for RSS?",
        "ques_type": 2,
        "options": [
            "Find the highest price and set it as the maximal threshold for RSS.",
            "Calculate the average house price and compute the difference.",
            "Find the lowest price and set it as the minimal threshold for RSS.",
            "Calculate the median house price and compute the difference."
        ],
        "score": "Calculate the average house price and compute the difference."
    },
    {
        "title": "You are applying k-means clustering on geospatial data of cities in the USA to determine cities falling in the same climate region. The available data contains latitudes and longitudes of cities.  Which distance metric should you use?",
        "ques_type": 2,
        "options": [
            "Actual distance between the cities using existing roads.",
            "Manhattan distance using latitude and longitude.",
            "Correlation coefficient using latitude and longitude.",
            "Euclidean distance using latitude and longitude."
        ],
        "score": "Euclidean distance using latitude and longitude."
    }
]
'@

# Remove old A2 content (the shared string moves up into A1)
$ws.Range("A2").ClearContents()

# A1 previously held a numeric 0 with bold/bordered/centered styling;
# strip that formatting back to the workbook default before writing
# the (now relocated) shared string into it.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $newText

# Writing a value that contains embedded newlines makes the engine
# auto-expand the row height (and pin customHeight); AutoFit() re-measures
# and clears that pinned/custom flag again, matching a freshly-saved row.
$ws.Rows.Item(1).AutoFit()
